$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: 2024-12-13T10:10:51-03:00 -> 2024-12-16T14:50:05-03:00
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true
$ws.Range("B17").Value = "'true"
